$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# (Target stored widths are 16.42578125 and 15.7109375; the host's column
# width model snaps to whole-pixel / 1/6-character increments, so we pick
# the ColumnWidth inputs whose rounded result is closest to the targets.)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.833333333333332

# Update cell values
$ws.Range("A1").Value = 0.012952780382252598
$ws.Range("B1").Value = -0.012952780721271543

$ws.Range("A2").Value = -0.047757753377730507
$ws.Range("B2").Value = 0.047757753063317601

$ws.Range("A3").Value = -0.00084126544920254029
$ws.Range("B3").Value = 0.00084126511792744516

$ws.Range("A4").Value = 0.025852928142408835
$ws.Range("B4").Value = -0.025852928462596125

$ws.Range("A5").Value = -0.00062626862750260795
$ws.Range("B5").Value = 0.00062626828076670288
